$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6959.5
$ws.Range("I40").Value = 13449.25
$ws.Range("J40").Value = 2633
$ws.Range("K40").Value = 13449.25
$ws.Range("L40").Value = 2633
$ws.Range("M40").Value = -13274.25
$ws.Range("N40").Value = -2983
$ws.Range("H58").Value = 13533.9375
$ws.Range("I58").Value = 1412.5
$ws.Range("J58").Value = 20806.8
$ws.Range("K58").Value = 4237.5
$ws.Range("L58").Value = 62420.39999999999
$ws.Range("M58").Value = -4087.5
$ws.Range("N58").Value = -62720.39999999999
$ws.Range("H92").Value = 522.3077
$ws.Range("I92").Value = 547.8889
$ws.Range("J92").Value = 464.75
$ws.Range("K92").Value = 547.8889
$ws.Range("L92").Value = 464.75
$ws.Range("M92").Value = 700.1111
$ws.Range("N92").Value = -2960.75
$ws.Range("H103").Value = 571.9655
$ws.Range("I103").Value = 556.8261
$ws.Range("J103").Value = 630
$ws.Range("K103").Value = 1670.4783
$ws.Range("L103").Value = 1890
$ws.Range("M103").Value = -1084.4783
$ws.Range("N103").Value = -3062
$ws.Range("H113").Value = 2850.5
$ws.Range("I113").Value = 2328
$ws.Range("J113").Value = 3212.2307
$ws.Range("K113").Value = 2328
$ws.Range("L113").Value = 3212.2307
$ws.Range("M113").Value = 926
$ws.Range("N113").Value = -9720.2307

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3858.4443
$ws.Range("I32").Value = 3858.4443
$ws.Range("K32").Value = 3858.4443
$ws.Range("M32").Value = -3571.4443
$ws.Range("H61").Value = 965.5
$ws.Range("I61").Value = 969.375
$ws.Range("J61").Value = 950
$ws.Range("K61").Value = 969.375
$ws.Range("L61").Value = 950
$ws.Range("M61").Value = -757.375
$ws.Range("N61").Value = -1374
$ws.Range("H111").Value = 88888
$ws.Range("J111").Value = 88888
$ws.Range("L111").Value = 88888
$ws.Range("N111").Value = -97068
$ws.Range("H132").Value = 1731.4286
$ws.Range("I132").Value = 1384
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 4152
$ws.Range("L132").Value = 7800
$ws.Range("M132").Value = -1622
$ws.Range("N132").Value = -12860
$ws.Range("H136").Value = 965.5
$ws.Range("I136").Value = 969.375
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 2908.125
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -358.125
$ws.Range("N136").Value = -7950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2260.7646
$ws.Range("I94").Value = 2161.4
$ws.Range("K94").Value = 2161.4
$ws.Range("M94").Value = -1710.4
$ws.Range("H97").Value = 14972.667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 14972.667
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 14972.667
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -16954.667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 16690.572
$ws.Range("I23").Value = 8997
$ws.Range("K23").Value = 8997
$ws.Range("M23").Value = -8757
$ws.Range("H26").Value = 10020
$ws.Range("I26").Value = 10020
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 10020
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -9733
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 16690.572
$ws.Range("I27").Value = 8997
$ws.Range("K27").Value = 8997
$ws.Range("M27").Value = -8805
$ws.Range("H105").Value = 3234.6365
$ws.Range("I105").Value = 2568.625
$ws.Range("J105").Value = 5010.6665
$ws.Range("K105").Value = 2568.625
$ws.Range("L105").Value = 5010.6665
$ws.Range("M105").Value = -821.625
$ws.Range("N105").Value = -8504.666499999999
$ws.Range("H130").Value = 51997.5
$ws.Range("J130").Value = 51997.5
$ws.Range("L130").Value = 51997.5
$ws.Range("N130").Value = -62037.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 880.36365
$ws.Range("I86").Value = 346
$ws.Range("J86").Value = 999.1111
$ws.Range("K86").Value = 1038
$ws.Range("L86").Value = 2997.3333
$ws.Range("M86").Value = 148
$ws.Range("N86").Value = -5369.3333
$ws.Range("H89").Value = 880.36365
$ws.Range("I89").Value = 346
$ws.Range("J89").Value = 999.1111
$ws.Range("K89").Value = 3114
$ws.Range("L89").Value = 8991.999899999999
$ws.Range("M89").Value = 2814
$ws.Range("N89").Value = -20847.9999
$ws.Range("H112").Value = 2695.6667
$ws.Range("H122").Value = 1317.091
$ws.Range("I122").Value = 741.2857
$ws.Range("J122").Value = 2324.75
$ws.Range("K122").Value = 6671.571300000001
$ws.Range("L122").Value = 20922.75
$ws.Range("M122").Value = -4221.571300000001
$ws.Range("N122").Value = -25822.75
$ws.Range("H129").Value = 2668.5715
$ws.Range("I129").Value = 830.6667
$ws.Range("J129").Value = 4047
$ws.Range("K129").Value = 2492.0001
$ws.Range("L129").Value = 12141
$ws.Range("M129").Value = 2507.9999
$ws.Range("N129").Value = -22141

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 30999.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30999.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30999.5
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -32191.5
$ws.Range("H97").Value = 978
$ws.Range("I97").Value = 725.8333
$ws.Range("J97").Value = 1482.3334
$ws.Range("K97").Value = 725.8333
$ws.Range("L97").Value = 1482.3334
$ws.Range("M97").Value = -229.8333
$ws.Range("N97").Value = -2474.3334
$ws.Range("H107").Value = 990.46155
$ws.Range("J107").Value = 873
$ws.Range("L107").Value = 873
$ws.Range("N107").Value = -4713
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920
$ws.Range("H132").Value = 2720.4614
$ws.Range("I132").Value = 2552.889
$ws.Range("K132").Value = 7658.667
$ws.Range("M132").Value = -5128.667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 90913590
$ws.Range("J136").Value = 250003860
$ws.Range("L136").Value = 750011580
$ws.Range("N136").Value = -750016680

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 4500
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H30").Value = 24381.2
$ws.Range("J30").Value = 24976.5
$ws.Range("L30").Value = 24976.5
$ws.Range("N30").Value = -25190.5
$ws.Range("H43").Value = 7749999.5
$ws.Range("I43").Value = 7749999.5
$ws.Range("K43").Value = 7749999.5
$ws.Range("M43").Value = -7749850.5
$ws.Range("H69").Value = 74999
$ws.Range("J69").Value = 74999
$ws.Range("L69").Value = 74999
$ws.Range("N69").Value = -76497
$ws.Range("H72").Value = 74999
$ws.Range("J72").Value = 74999
$ws.Range("L72").Value = 224997
$ws.Range("N72").Value = -232485
$ws.Range("H107").Value = 3423.8
$ws.Range("I107").Value = 1055
$ws.Range("J107").Value = 5003
$ws.Range("K107").Value = 3165
$ws.Range("L107").Value = 15009
$ws.Range("M107").Value = -1245
$ws.Range("N107").Value = -18849
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = 0
$ws.Range("H122").Value = 3632.2068
$ws.Range("I122").Value = 3288.5
$ws.Range("K122").Value = 9865.5
$ws.Range("M122").Value = -7415.5

Write-Output "Applied all Alpha_Profits cell updates."